# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K") values are recalculated for each row (2-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..33 (computed from the regenerated stat, replacing
# the old Strike# values that used to live in column G).
$kValues = @{
    2  = 1
    3  = 2
    4  = 4
    5  = 0
    6  = 5
    7  = 0
    8  = 4
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 4
    15 = 1
    16 = 6
    17 = 2
    18 = 5
    19 = 2
    20 = 6
    21 = 3
    22 = 9
    23 = 5
    24 = 4
    25 = 3
    26 = 0
    27 = 4
    28 = 4
    29 = 3
    30 = 4
    31 = 2
    32 = 3
    33 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
